# Baza_danych.xlsx — update the "Mieszkania" listing table.
#
# The Gdansk-area listings that used to occupy rows 2-5 are replaced by the
# Lublin / Gliwice / Katowice listings that used to live in rows 6-8 (each
# repeated twice, filling rows 2-4 and 5-7), and the now-redundant row 8 is
# removed, shrinking the used range from A1:O8 to A1:O7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the last row first so the sheet's used range becomes A1:O7.
$ws.Rows("8:8").Delete()

# The "liczba_pokoi" (D), "pietro" (E) and "rok_budowy" (G) columns hold
# plain-looking numbers (e.g. "2", "1960") but the source data stores them
# as text, matching the rest of the sheet (no numeric formatting anywhere
# in this workbook). Force those columns to Text before writing so the
# digit strings round-trip as text instead of being auto-coerced to numbers.
$ws.Range("D2:E7").NumberFormat = "@"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G5").NumberFormat = "@"

# --- Row 2: 610 000 zl / Lublin, Rury, Konrada Wallenroda ---
$ws.Range("A2").Value = "610 000 zł"
$ws.Range("B2").Value = "10854 zł/m²"
$ws.Range("C2").Value = "56,20 m²"
$ws.Range("D2").Value = "2"
$ws.Range("E2").Value = "2"
$ws.Range("F2").Value = "wtórny"
$ws.Range("G2").Value = "1960"
$ws.Range("H2").ClearContents()
$ws.Range("I2").Value = "Lubelskie"
$ws.Range("J2").Value = "Lublin"
$ws.Range("K2").Value = "Lublin"
$ws.Range("L2").Value = "Lublin"
$ws.Range("M2").Value = "Rury"
$ws.Range("N2").Value = "Konrada Wallenroda"
$ws.Range("O2").Value = "https://www.otodom.pl/pl/oferta/mieszkanie-2-pokojowe-rury-po-remoncie-ID4xE5n"

# --- Row 3: 239 000 zl / Gliwice, Zatorze ---
$ws.Range("A3").Value = "239 000 zł"
$ws.Range("B3").Value = "9192 zł/m²"
$ws.Range("C3").Value = "26 m²"
$ws.Range("D3").Value = "1"
$ws.Range("E3").Value = "3"
$ws.Range("F3").Value = "wtórny"
$ws.Range("G3").ClearContents()
$ws.Range("H3").ClearContents()
$ws.Range("I3").Value = "Śląskie"
$ws.Range("J3").Value = "Gliwice"
$ws.Range("K3").Value = "Gliwice"
$ws.Range("L3").Value = "Gliwice"
$ws.Range("M3").Value = "Zatorze"
$ws.Range("N3").ClearContents()
$ws.Range("O3").Value = "https://www.otodom.pl/pl/oferta/przytulna-kawalerka-z-klimatyzacja-idealna-do-zamieszkania-ID4v2eg"

# --- Row 4: 465 430 zl / Katowice, Welnowiec-Jozefowiec, Tomasza Kotlarza ---
$ws.Range("A4").Value = "465 430 zł"
$ws.Range("B4").Value = "10900 zł/m²"
$ws.Range("C4").Value = "42,70 m²"
$ws.Range("D4").Value = "3"
$ws.Range("E4").Value = "1"
$ws.Range("F4").Value = "wtórny"
$ws.Range("G4").ClearContents()
$ws.Range("H4").Value = "wielka płyta"
$ws.Range("I4").Value = "Śląskie"
$ws.Range("J4").Value = "Katowice"
$ws.Range("K4").Value = "Katowice"
$ws.Range("L4").Value = "Katowice"
$ws.Range("M4").Value = "Wełnowiec-Józefowiec"
$ws.Range("N4").Value = "Tomasza Kotlarza"
$ws.Range("O4").Value = "https://www.otodom.pl/pl/oferta/odkryj-stylowe-wnetrze-tuz-obok-parku-slaskiego-i-silesii-ID4xEYt"

# --- Row 5: 610 000 zl / Lublin, Rury, Konrada Wallenroda (2nd copy) ---
$ws.Range("A5").Value = "610 000 zł"
$ws.Range("B5").Value = "10854 zł/m²"
$ws.Range("C5").Value = "56,20 m²"
$ws.Range("D5").Value = "2"
$ws.Range("E5").Value = "2"
$ws.Range("F5").Value = "wtórny"
$ws.Range("G5").Value = "1960"
$ws.Range("H5").ClearContents()
$ws.Range("I5").Value = "Lubelskie"
$ws.Range("J5").Value = "Lublin"
$ws.Range("K5").Value = "Lublin"
$ws.Range("L5").Value = "Lublin"
$ws.Range("M5").Value = "Rury"
$ws.Range("N5").Value = "Konrada Wallenroda"
$ws.Range("O5").Value = "https://www.otodom.pl/pl/oferta/mieszkanie-2-pokojowe-rury-po-remoncie-ID4xE5n"

# --- Row 6: 239 000 zl / Gliwice, Zatorze (2nd copy) ---
$ws.Range("A6").Value = "239 000 zł"
$ws.Range("B6").Value = "9192 zł/m²"
$ws.Range("C6").Value = "26 m²"
$ws.Range("D6").Value = "1"
$ws.Range("E6").Value = "3"
$ws.Range("F6").Value = "wtórny"
$ws.Range("G6").ClearContents()
$ws.Range("H6").ClearContents()
$ws.Range("I6").Value = "Śląskie"
$ws.Range("J6").Value = "Gliwice"
$ws.Range("K6").Value = "Gliwice"
$ws.Range("L6").Value = "Gliwice"
$ws.Range("M6").Value = "Zatorze"
$ws.Range("N6").ClearContents()
$ws.Range("O6").Value = "https://www.otodom.pl/pl/oferta/przytulna-kawalerka-z-klimatyzacja-idealna-do-zamieszkania-ID4v2eg"

# --- Row 7: 465 430 zl / Katowice, Welnowiec-Jozefowiec, Tomasza Kotlarza (2nd copy) ---
$ws.Range("A7").Value = "465 430 zł"
$ws.Range("B7").Value = "10900 zł/m²"
$ws.Range("C7").Value = "42,70 m²"
$ws.Range("D7").Value = "3"
$ws.Range("E7").Value = "1"
$ws.Range("F7").Value = "wtórny"
$ws.Range("G7").ClearContents()
$ws.Range("H7").Value = "wielka płyta"
$ws.Range("I7").Value = "Śląskie"
$ws.Range("J7").Value = "Katowice"
$ws.Range("K7").Value = "Katowice"
$ws.Range("L7").Value = "Katowice"
$ws.Range("M7").Value = "Wełnowiec-Józefowiec"
$ws.Range("N7").Value = "Tomasza Kotlarza"
$ws.Range("O7").Value = "https://www.otodom.pl/pl/oferta/odkryj-stylowe-wnetrze-tuz-obok-parku-slaskiego-i-silesii-ID4xEYt"
